$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Column D ("reps") for WEEK 2, WEEK 3 and WEEK 4 blocks,
# matching / extending the scheme already present for WEEK 1 (rows 2-4).

# WEEK 2 block (rows 8-10) - same scheme as WEEK 1
$ws.Range("D8").Value = "1x20, 5x1"
$ws.Range("D9").Value = "1x20, 10x1"
$ws.Range("D10").Value = "1x20, 1x10, 1x2.5"

# WEEK 3 block (rows 14-16)
$ws.Range("D14").Value = "1x20, 5x1, 1x2.5"
$ws.Range("D15").Value = "1x20, 1x10, 1x2.5"
$ws.Range("D16").Value = "1x20, 1x10, 1x5, 1x2.5"

# WEEK 4 block (rows 20-22)
$ws.Range("D20").Value = "1x20"
$ws.Range("D21").Value = "1x20"
$ws.Range("D22").Value = "1x20"

# Update the view: move the active selection to D22 (last entry just
# filled in) and scroll the window down so row 11 is the top visible row.
$ws.Range("D22").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
